# Update "想去人数" (column F) values for rows 3-13 on the
# "展览" and "全部类型" worksheets, which hold identical data.

$wb = $excel.ActiveWorkbook

# Mapping of row -> new value for column F
$updates = @{
    3  = 1389
    4  = 6707
    5  = 363
    6  = 194
    7  = 3284
    8  = 17
    9  = 12
    10 = 42
    11 = 816
    12 = 250
    13 = 5401
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
